$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Remove the slides that are no longer part of the deck.
#    The final deck keeps only the former slides 1, 7, 8 and 11 (in that
#    order), which will become the new slides 1-4.
#    Delete from the highest index down so earlier indices stay valid.
# ---------------------------------------------------------------------------
$slidesToDelete = @(17,16,15,14,13,12,10,9,6,5,4,3,2)
foreach ($idx in $slidesToDelete) {
    $p.Slides.Item($idx).Delete()
}

# ---------------------------------------------------------------------------
# 2) Strip the bold emphasis that was consolidated away on the remaining
#    slides (now indices 2, 3 and 4).
# ---------------------------------------------------------------------------

# -- New slide 2 (was "Slide 7: Why This Solution?") ------------------------
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Cell(4,1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl2.Cell(4,2).Shape.TextFrame.TextRange.Font.Bold = $false

# -- New slide 3 (was "Slide 8: Business Value - Financial Impact") --------
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
$tbl3.Cell(1,1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(1,2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(6,1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(6,2).Shape.TextFrame.TextRange.Font.Bold = $false

# -- New slide 4 (was "Slide 11: Risk Mitigation") --------------------------
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
$tbl4.Cell(1,1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(1,2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(1,3).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4,1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4,2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4,3).Shape.TextFrame.TextRange.Font.Bold = $false

Write-Host "Final slide count: $($p.Slides.Count)"
